$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, [string]$val)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

# Row -> (D value or $null, E value)
$updates = @(
    @{ Row = 2;  D = "70.848.08";  E = "  +5.56%  " },
    @{ Row = 3;  D = "3.655.35";   E = "  +5.40%  " },
    @{ Row = 4;  D = $null;        E = "  -0.11%  " },
    @{ Row = 5;  D = "593.41";     E = "  +1.27%  " },
    @{ Row = 6;  D = "195.15";     E = "  +3.56%  " },
    @{ Row = 7;  D = "0.650";      E = "  +2.91%  " },
    @{ Row = 8;  D = "3.649.09";   E = "  +5.38%  " },
    @{ Row = 9;  D = $null;        E = "  -0.07%  " },
    @{ Row = 10; D = $null;        E = "  +5.24%  " },
    @{ Row = 11; D = "0.678";      E = "  +4.51%  " },
    @{ Row = 12; D = "58.67";      E = "  +3.18%  " },
    @{ Row = 13; D = $null;        E = "  +5.37%  " },
    @{ Row = 14; D = "9.99";       E = "  +5.86%  " },
    @{ Row = 15; D = "4.234.95";   E = "  +5.13%  " },
    @{ Row = 16; D = "19.97";      E = "  +6.41%  " },
    @{ Row = 17; D = "3.653.12";   E = "  +5.23%  " },
    @{ Row = 18; D = "70.808.65";  E = "  +5.41%  " },
    @{ Row = 19; D = "12.83";      E = "  +5.52%  " },
    @{ Row = 20; D = $null;        E = "  +2.45%  " },
    @{ Row = 21; D = $null;        E = "  +5.30%  " },
    @{ Row = 22; D = "492.30";     E = "  +1.30%  " },
    @{ Row = 23; D = "18.94";      E = "  +12.59%  " },
    @{ Row = 24; D = "5.36";       E = "  -0.06%  " },
    @{ Row = 25; D = $null;        E = "  +0.86%  " },
    @{ Row = 26; D = "91.96";      E = "  +2.58%  " },
    @{ Row = 27; D = "3.18";       E = "  +8.11%  " },
    @{ Row = 28; D = $null;        E = "  +5.08%  " },
    @{ Row = 29; D = "9.64";       E = "  +5.81%  " },
    @{ Row = 30; D = "33.14";      E = "  +5.59%  " },
    @{ Row = 31; D = "7.91";       E = "  +10.17%  " },
    @{ Row = 32; D = "0.122";      E = "  +8.90%  " },
    @{ Row = 33; D = "631.73";     E = "  +5.86%  " },
    @{ Row = 34; D = "12.34";      E = "  +5.08%  " },
    @{ Row = 35; D = "65.77";      E = "  +2.62%  " },
    @{ Row = 36; D = "40.72";      E = "  +10.97%  " },
    @{ Row = 37; D = "0.0₃0841";   E = "  +10.97%  " },
    @{ Row = 38; D = $null;        E = "  +7.78%  " },
    @{ Row = 39; D = $null;        E = "  -1.20%  " },
    @{ Row = 40; D = $null;        E = "  -0.01%  " },
    @{ Row = 41; D = $null;        E = "  +1.65%  " },
    @{ Row = 42; D = "3.321.93";   E = "  +2.64%  " },
    @{ Row = 44; D = $null;        E = "  +14.25%  " },
    @{ Row = 45; D = $null;        E = "  +5.71%  " },
    @{ Row = 46; D = "2.95";       E = "  +5.54%  " },
    @{ Row = 47; D = "3.31";       E = "  +1.06%  " },
    @{ Row = 48; D = $null;        E = "  +2.82%  " },
    @{ Row = 49; D = $null;        E = "  +5.96%  " },
    @{ Row = 50; D = $null;        E = "  +1.90%  " },
    @{ Row = 51; D = "1.00";       E = "  -0.12%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        Set-TextCell $ws.Cells.Item($r, 4) $u.D
    }
    Set-TextCell $ws.Cells.Item($r, 5) $u.E
}
